{"js": "const replacements = [\n  [\"2024-11-23 Saturday\", \"2024-11-24 Sunday\"],\n  [\"82\u00d774=\", \"57\u00d779=\"],\n  [\"43\u00d738=\", \"73\u00d750=\"],\n  [\"54\u00d734=\", \"42\u00d781=\"],\n  [\"86\u00d793=\", \"85\u00d712=\"],\n  [\"85\u00d750=\", \"24\u00d724=\"],\n  [\"35\u00d750=\", \"81\u00d761=\"],\n  [\"86\u00d714=\", \"24\u00d746=\"],\n  [\"95\u00d755=\", \"48\u00d734=\"],\n  [\"25\u00d767=\", \"38\u00d729=\"],\n  [\"69\u00d735=\", \"57\u00d787=\"],\n  [\"12\u00d739=\", \"22\u00d772=\"],\n  [\"28\u00d790=\", \"60\u00d751=\"],\n  [\"70\u00d749=\", \"70\u00d779=\"],\n  [\"15\u00d754=\", \"53\u00d724=\"],\n  [\"21\u00d756=\", \"55\u00d773=\"],\n  [\"18\u00d792=\", \"80\u00d751=\"],\n  [\"25\u00d796=\", \"77\u00d767=\"],\n  [\"31\u00d730=\", \"61\u00d795=\"],\n  [\"47\u00d762=\", \"72\u00d711=\"],\n  [\"84\u00d765=\", \"57\u00d718=\"],\n  [\"76\u00d782=\", \"56\u00d755=\"],\n  [\"33\u00d755=\", \"28\u00d786=\"],\n  [\"26\u00d719=\", \"92\u00d769=\"],\n  [\"41\u00d738=\", \"41\u00d722=\"],\n  [\"48\u00d720=\", \"99\u00d792=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"2024-11-23 Saturday\"\n$find.Replacement.Text = \"2024-11-24 Sunday\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"82\u00d774=\"\n$find.Replacement.Text = \"57\u00d779=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"43\u00d738=\"\n$find.Replacement.Text = \"73\u00d750=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"54\u00d734=\"\n$find.Replacement.Text = \"42\u00d781=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"86\u00d793=\"\n$find.Replacement.Text = \"85\u00d712=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"85\u00d750=\"\n$find.Replacement.Text = \"24\u00d724=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"35\u00d750=\"\n$find.Replacement.Text = \"81\u00d761=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"86\u00d714=\"\n$find.Replacement.Text = \"24\u00d746=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"95\u00d755=\"\n$find.Replacement.Text = \"48\u00d734=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"25\u00d767=\"\n$find.Replacement.Text = \"38\u00d729=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"69\u00d735=\"\n$find.Replacement.Text = \"57\u00d787=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"12\u00d739=\"\n$find.Replacement.Text = \"22\u00d772=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"28\u00d790=\"\n$find.Replacement.Text = \"60\u00d751=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"70\u00d749=\"\n$find.Replacement.Text = \"70\u00d779=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"15\u00d754=\"\n$find.Replacement.Text = \"53\u00d724=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"21\u00d756=\"\n$find.Replacement.Text = \"55\u00d773=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"18\u00d792=\"\n$find.Replacement.Text = \"80\u00d751=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"25\u00d796=\"\n$find.Replacement.Text = \"77\u00d767=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"31\u00d730=\"\n$find.Replacement.Text = \"61\u00d795=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"47\u00d762=\"\n$find.Replacement.Text = \"72\u00d711=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"84\u00d765=\"\n$find.Replacement.Text = \"57\u00d718=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"76\u00d782=\"\n$find.Replacement.Text = \"56\u00d755=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"33\u00d755=\"\n$find.Replacement.Text = \"28\u00d786=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"26\u00d719=\"\n$find.Replacement.Text = \"92\u00d769=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"41\u00d738=\"\n$find.Replacement.Text = \"41\u00d722=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"48\u00d720=\"\n$find.Replacement.Text = \"99\u00d792=\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
